$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.873501
$ws.Range("H2").Value = 2.620503
$ws.Range("I2").Value = 0.1166943280075418
$ws.Range("J2").Value = 0.1166943280075418
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.873501
$ws.Range("N2").Value = 2.620503
$ws.Range("O2").Value = 0.1166943280075418
$ws.Range("P2").Value = 0.1166943280075418
$ws.Range("Q2").Value = 0.763003997001
$ws.Range("R2").Value = 6.867035973008999
$ws.Range("S2").Value = 0.01361756618913175
$ws.Range("T2").Value = 0.01361756618913175

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.873501
$ws.Range("H3").Value = 2.620503
$ws.Range("I3").Value = 0.1166943280075418
$ws.Range("J3").Value = 0.1166943280075418
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.327332
$ws.Range("N3").Value = 0.9819960000000001
$ws.Range("O3").Value = 0.04372952953158002
$ws.Range("P3").Value = 0.04372952953158002
$ws.Range("Q3").Value = 0.285924829332
$ws.Range("R3").Value = 2.573323463988
$ws.Range("S3").Value = 0.005102988062773685
$ws.Range("T3").Value = 0.005102988062773685

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.873501
$ws.Range("H4").Value = 2.620503
$ws.Range("I4").Value = 0.1166943280075418
$ws.Range("J4").Value = 0.1166943280075418
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.284543666666667
$ws.Range("N4").Value = 18.853631
$ws.Range("O4").Value = 0.8395761424608782
$ws.Range("P4").Value = 0.8395761424608781
$ws.Range("Q4").Value = 5.489555177377
$ws.Range("R4").Value = 49.405996596393
$ws.Range("S4").Value = 0.09797377375563636
$ws.Range("T4").Value = 0.09797377375563635

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.327332
$ws.Range("H5").Value = 0.9819960000000001
$ws.Range("I5").Value = 0.04372952953158002
$ws.Range("J5").Value = 0.04372952953158002
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.873501
$ws.Range("N5").Value = 2.620503
$ws.Range("O5").Value = 0.1166943280075418
$ws.Range("P5").Value = 0.1166943280075418
$ws.Range("Q5").Value = 0.285924829332
$ws.Range("R5").Value = 2.573323463988
$ws.Range("S5").Value = 0.005102988062773685
$ws.Range("T5").Value = 0.005102988062773685

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.327332
$ws.Range("H6").Value = 0.9819960000000001
$ws.Range("I6").Value = 0.04372952953158002
$ws.Range("J6").Value = 0.04372952953158002
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.327332
$ws.Range("N6").Value = 0.9819960000000001
$ws.Range("O6").Value = 0.04372952953158002
$ws.Range("P6").Value = 0.04372952953158002
$ws.Range("Q6").Value = 0.107146238224
$ws.Range("R6").Value = 0.9643161440160002
$ws.Range("S6").Value = 0.001912271753053329
$ws.Range("T6").Value = 0.001912271753053329

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.327332
$ws.Range("H7").Value = 0.9819960000000001
$ws.Range("I7").Value = 0.04372952953158002
$ws.Range("J7").Value = 0.04372952953158002
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 6.284543666666667
$ws.Range("N7").Value = 18.853631
$ws.Range("O7").Value = 0.8395761424608782
$ws.Range("P7").Value = 0.8395761424608781
$ws.Range("Q7").Value = 2.057132247497333
$ws.Range("R7").Value = 18.514190227476
$ws.Range("S7").Value = 0.03671426971575301
$ws.Range("T7").Value = 0.036714269715753

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 6.284543666666667
$ws.Range("H8").Value = 18.853631
$ws.Range("I8").Value = 0.8395761424608782
$ws.Range("J8").Value = 0.8395761424608781
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.873501
$ws.Range("N8").Value = 2.620503
$ws.Range("O8").Value = 0.1166943280075418
$ws.Range("P8").Value = 0.1166943280075418
$ws.Range("Q8").Value = 5.489555177377
$ws.Range("R8").Value = 49.405996596393
$ws.Range("S8").Value = 0.09797377375563636
$ws.Range("T8").Value = 0.09797377375563635

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 6.284543666666667
$ws.Range("H9").Value = 18.853631
$ws.Range("I9").Value = 0.8395761424608782
$ws.Range("J9").Value = 0.8395761424608781
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.327332
$ws.Range("N9").Value = 0.9819960000000001
$ws.Range("O9").Value = 0.04372952953158002
$ws.Range("P9").Value = 0.04372952953158002
$ws.Range("Q9").Value = 2.057132247497333
$ws.Range("R9").Value = 18.514190227476
$ws.Range("S9").Value = 0.03671426971575301
$ws.Range("T9").Value = 0.036714269715753

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 6.284543666666667
$ws.Range("H10").Value = 18.853631
$ws.Range("I10").Value = 0.8395761424608782
$ws.Range("J10").Value = 0.8395761424608781
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 6.284543666666667
$ws.Range("N10").Value = 18.853631
$ws.Range("O10").Value = 0.8395761424608782
$ws.Range("P10").Value = 0.8395761424608781
$ws.Range("Q10").Value = 39.49548909824011
$ws.Range("R10").Value = 355.459401884161
$ws.Range("S10").Value = 0.7048880989894888
$ws.Range("T10").Value = 0.7048880989894887
